$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Orthography fix: replace the digit "1" (used as a stand-in for the
# palochka letter) with the Latin capital letter "I" throughout the
# Botlikh-language column (column B), per Azaev (1974).

$ws.Range("B2").Value = "Керпеч букъи буккала, гьечIи ссеъа унщи хварде, ехха бехутIу лълъенила тIинду биччи бигъай амал игье."
$ws.Range("B3").Value = "Биччи бигъаата-кIо ригьуди бакьа никкула бигье."
$ws.Range("B4").Value = "Эхха гьеъала рулIу биччи цIекIабалъиди мерхьу булIа-риуди тIоргуда биччи данда игье."
$ws.Range("B6").Value = "Керпеч букъихо адамий кумак игье цевла гьекIващуди."
$ws.Range("B7").Value = "Гьув гьекIващуди биччи бекьара барлIе."
$ws.Range("B8").Value = "Керпеч букъихо адамиди кепилъи гьину биччила букIо реъабалъиди биччи тIибдай амал игье."
$ws.Range("B9").Value = "Эхха кепла гьирцIо гъижихи букIай амале. "
$ws.Range("B10").Value = "КIейхоб букхата-кIо ригьуди чинкурла ибху кеп гьинуку аржай букке."
$ws.Range("B12").Value = "Букъу булъариуди керпеч букъуй бете, эхха уда — ара ссардай амал игье бекъа бекъуI-талу."

# Reflect the last-edited cell as the active selection, matching the
# author's editing session ending on B8.
$ws.Range("B8").Select()
